$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "BBB"
$ws.Range("D9").Value = "A"
$ws.Range("D10").Value = "A"
$ws.Range("D11").Value = "A"
$ws.Range("D12").Value = "A"
$ws.Range("D33").Value = "BBB"
$ws.Range("D39").Value = "B"
$ws.Range("D51").Value = "BB"
$ws.Range("D62").Value = "A"
$ws.Range("D76").Value = "BBB"
$ws.Range("D85").Value = "B"
$ws.Range("D87").Value = "BBB"
$ws.Range("D95").Value = "A"
$ws.Range("D113").Value = "A"
$ws.Range("D116").Value = "A"
$ws.Range("D118").Value = "BBB"
$ws.Range("D140").Value = "BBB"
$ws.Range("D141").Value = "BB"
$ws.Range("D143").Value = "AA"
$ws.Range("D146").Value = "A"
$ws.Range("D156").Value = "B"
$ws.Range("D157").Value = "B"
$ws.Range("D158").Value = "BB"
$ws.Range("D159").Value = "A"
$ws.Range("D161").Value = "A"
$ws.Range("D167").Value = "A"
$ws.Range("D169").Value = "BBB"
$ws.Range("D174").Value = "BB"
$ws.Range("D180").Value = "A"
$ws.Range("D182").Value = "A"
$ws.Range("D184").Value = "B"
$ws.Range("D187").Value = "AAA"
$ws.Range("D195").Value = "BBB"
$ws.Range("D212").Value = "B"
$ws.Range("D226").Value = "BBB"
$ws.Range("D229").Value = "B"
$ws.Range("D233").Value = "BBB"
$ws.Range("D239").Value = "A"
$ws.Range("D249").Value = "B"
$ws.Range("D255").Value = "BBB"
$ws.Range("D270").Value = "BB"
$ws.Range("D273").Value = "BBB"
$ws.Range("D283").Value = "BB"
$ws.Range("D286").Value = "BBB"
$ws.Range("D288").Value = "B"
$ws.Range("D290").Value = "BBB"
$ws.Range("D291").Value = "BBB"
$ws.Range("D300").Value = "A"
$ws.Range("D301").Value = "A"
$ws.Range("D302").Value = "A"
$ws.Range("D303").Value = "BB"
$ws.Range("D311").Value = "BBB"
$ws.Range("D313").Value = "B"
$ws.Range("D342").Value = "BBB"
$ws.Range("D349").Value = "B"
$ws.Range("D360").Value = "AA"
$ws.Range("D361").Value = "A"
$ws.Range("D370").Value = "BBB"
$ws.Range("D385").Value = "BB"
$ws.Range("D387").Value = "BBB"
$ws.Range("D394").Value = "BB"
